$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark F9 with "X" (new entry)
$ws.Range("F9").Value = "X"

# Clear F14 (the "X" there is removed; B14 keeps its "X")
$ws.Range("F14").ClearContents()

# Update the active selection to H16, matching the saved view state
$ws.Range("H16").Select()
